# Automatic update by jenkins
#
# Adds a "lat"/"lng" pair of geo-coordinate columns to the
# root_hospital_cities sheet, documents the two new attributes on the
# "attributes" metadata sheet, and leaves the workbook's active sheet /
# selection pointing at root_hospital_cities (tab 1) instead of entities
# (tab 5).

$wb = $excel.ActiveWorkbook

# --- sheet "root_hospital_cities": add lat/lng columns -------------------
$wsCities = $wb.Worksheets.Item("root_hospital_cities")

$wsCities.Range("B1").Value = "lat"
$wsCities.Range("C1").Value = "lng"
$wsCities.Range("B1:C1").Font.Color = 0

$wsCities.Range("B2").Value = 40.712784
$wsCities.Range("C2").Value = -74.005941
$wsCities.Range("B3").Value = 37.151165
$wsCities.Range("C3").Value = -88.731998
$wsCities.Range("B2:C3").Font.Color = 0
$wsCities.Range("B2:C3").NumberFormat = "0.000000"

# --- sheet "attributes": document the new lat/lng attributes -------------
$wsAttributes = $wb.Worksheets.Item("attributes")

$wsAttributes.Rows("3:4").Insert()

$wsAttributes.Range("A3").Value = "lat"
$wsAttributes.Range("B3").Value = "root_hospital_cities"
$wsAttributes.Range("A4").Value = "lng"
$wsAttributes.Range("B4").Value = "root_hospital_cities"

$wsAttributes.Range("G3").Value = "latitude in degrees"
$wsAttributes.Range("G4").Value = "longitude in degrees"

$wsAttributes.Range("C3").Value = "decimal"
$wsAttributes.Range("C4").Value = "decimal"

# the attribute sheet's remembered selection moves from the (now shifted)
# B5 cell to C5
$wsAttributes.Range("C5").Select()

# --- move the active tab / selection back to root_hospital_cities --------
$wsCities.Range("A1").Select()
